$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Cells.Item(2,1).Value = "Última actualización: 20:47:47"
$ws1.Cells.Item(3,1).Value = "Total filas: 499"
$ws1.Cells.Item(190,1).Value = "10:59:49"
$ws1.Cells.Item(190,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(190,4).Value = 82
$ws1.Cells.Item(191,1).Value = "12:21:08"
$ws1.Cells.Item(191,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(191,4).Value = 0
$ws1.Cells.Item(192,3).Value = "215A_EL PATO"
$ws1.Cells.Item(238,1).Value = "12:59:47"
$ws1.Cells.Item(238,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(238,4).Value = 51
$ws1.Cells.Item(239,1).Value = "11:56:55"
$ws1.Cells.Item(239,3).Value = "215A_EL PATO"
$ws1.Cells.Item(239,4).Value = 114
$ws1.Cells.Item(300,1).Value = "15:53:28"
$ws1.Cells.Item(300,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(300,4).Value = 3
$ws1.Cells.Item(301,1).Value = "14:24:16"
$ws1.Cells.Item(301,3).Value = "17_ROMERO"
$ws1.Cells.Item(301,4).Value = 92
$ws1.Cells.Item(332,1).Value = "16:31:51"
$ws1.Cells.Item(332,3).Value = "10_OLMOS"
$ws1.Cells.Item(332,4).Value = 25
$ws1.Cells.Item(333,1).Value = "15:22:17"
$ws1.Cells.Item(333,3).Value = "17_179 Y 38"
$ws1.Cells.Item(333,4).Value = 94
$ws1.Cells.Item(364,1).Value = "17:38:19"
$ws1.Cells.Item(364,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(364,4).Value = 14
$ws1.Cells.Item(365,1).Value = "15:53:28"
$ws1.Cells.Item(365,3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(365,4).Value = 119
$ws1.Cells.Item(366,1).Value = "17:51:34"
$ws1.Cells.Item(366,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(366,4).Value = 1
$ws1.Cells.Item(410,1).Value = "17:51:34"
$ws1.Cells.Item(410,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(410,4).Value = 85
$ws1.Cells.Item(411,1).Value = "18:33:37"
$ws1.Cells.Item(411,3).Value = "15_ABASTO"
$ws1.Cells.Item(411,4).Value = 43
$ws1.Cells.Item(457,3).Value = "215A_EL PATO"
$ws1.Cells.Item(458,3).Value = "14_ABASTO"
$ws1.Cells.Item(470,1).Value = "20:47:47"
$ws1.Cells.Item(470,2).Value = "20:48"
$ws1.Cells.Item(470,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(470,4).Value = 1
$ws1.Cells.Item(472,1).Value = "19:13:07"
$ws1.Cells.Item(472,2).Value = "20:52"
$ws1.Cells.Item(472,3).Value = "15_ABASTO"
$ws1.Cells.Item(472,4).Value = 99
$ws1.Cells.Item(473,1).Value = "20:47:47"
$ws1.Cells.Item(473,2).Value = "20:53"
$ws1.Cells.Item(473,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(473,4).Value = 6
$ws1.Cells.Item(474,1).Value = "20:33:52"
$ws1.Cells.Item(474,2).Value = "20:54"
$ws1.Cells.Item(474,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(474,4).Value = 21
$ws1.Cells.Item(475,1).Value = "20:12:40"
$ws1.Cells.Item(475,2).Value = "20:56"
$ws1.Cells.Item(475,4).Value = 44
$ws1.Cells.Item(476,1).Value = "19:52:18"
$ws1.Cells.Item(476,2).Value = "20:56"
$ws1.Cells.Item(476,3).Value = "10_OLMOS"
$ws1.Cells.Item(476,4).Value = 64
$ws1.Cells.Item(477,1).Value = "19:13:07"
$ws1.Cells.Item(477,2).Value = "20:57"
$ws1.Cells.Item(477,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(477,4).Value = 104
$ws1.Cells.Item(478,1).Value = "19:13:07"
$ws1.Cells.Item(478,2).Value = "21:04"
$ws1.Cells.Item(478,3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(478,4).Value = 111
$ws1.Cells.Item(479,1).Value = "20:12:40"
$ws1.Cells.Item(479,2).Value = "21:07"
$ws1.Cells.Item(479,4).Value = 55
$ws1.Cells.Item(480,1).Value = "20:33:52"
$ws1.Cells.Item(480,2).Value = "21:07"
$ws1.Cells.Item(480,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(480,4).Value = 34
$ws1.Cells.Item(481,1).Value = "19:13:07"
$ws1.Cells.Item(481,2).Value = "21:08"
$ws1.Cells.Item(481,3).Value = "215B_EL PATO"
$ws1.Cells.Item(481,4).Value = 115
$ws1.Cells.Item(482,2).Value = "21:08"
$ws1.Cells.Item(482,3).Value = "16_P MOR-167 Y 521"
$ws1.Cells.Item(482,4).Value = 56
$ws1.Cells.Item(483,2).Value = "21:16"
$ws1.Cells.Item(483,3).Value = "14_ABASTO"
$ws1.Cells.Item(483,4).Value = 97
$ws1.Cells.Item(484,1).Value = "20:12:40"
$ws1.Cells.Item(484,2).Value = "21:20"
$ws1.Cells.Item(484,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(484,4).Value = 68
$ws1.Cells.Item(485,2).Value = "21:21"
$ws1.Cells.Item(485,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(485,4).Value = 102
$ws1.Cells.Item(486,2).Value = "21:22"
$ws1.Cells.Item(486,3).Value = "15_ABASTO"
$ws1.Cells.Item(486,4).Value = 49
$ws1.Cells.Item(487,1).Value = "19:39:04"
$ws1.Cells.Item(487,2).Value = "21:23"
$ws1.Cells.Item(487,3).Value = "10_OLMOS"
$ws1.Cells.Item(487,4).Value = 104
$ws1.Cells.Item(488,1).Value = "20:33:52"
$ws1.Cells.Item(488,2).Value = "21:32"
$ws1.Cells.Item(488,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(488,4).Value = 59
$ws1.Cells.Item(489,1).Value = "20:33:52"
$ws1.Cells.Item(489,2).Value = "21:37"
$ws1.Cells.Item(489,3).Value = "14_ABASTO"
$ws1.Cells.Item(489,4).Value = 64
$ws1.Cells.Item(490,1).Value = "20:12:40"
$ws1.Cells.Item(490,2).Value = "21:37"
$ws1.Cells.Item(490,3).Value = "17_ROMERO"
$ws1.Cells.Item(490,4).Value = 85
$ws1.Cells.Item(491,1).Value = "19:52:18"
$ws1.Cells.Item(491,2).Value = "21:38"
$ws1.Cells.Item(491,3).Value = "14_ABASTO"
$ws1.Cells.Item(491,4).Value = 106
$ws1.Cells.Item(492,1).Value = "19:39:04"
$ws1.Cells.Item(492,2).Value = "21:38"
$ws1.Cells.Item(492,3).Value = "17_ROMERO"
$ws1.Cells.Item(492,4).Value = 119
$ws1.Cells.Item(493,1).Value = "20:33:52"
$ws1.Cells.Item(493,2).Value = "21:46"
$ws1.Cells.Item(493,3).Value = "215A_EL PATO"
$ws1.Cells.Item(493,4).Value = 73
$ws1.Cells.Item(494,1).Value = "19:52:18"
$ws1.Cells.Item(494,2).Value = "21:47"
$ws1.Cells.Item(494,3).Value = "215A_EL PATO"
$ws1.Cells.Item(495,1).Value = "20:12:40"
$ws1.Cells.Item(495,2).Value = "21:52"
$ws1.Cells.Item(495,3).Value = "10_OLMOS"
$ws1.Cells.Item(495,4).Value = 100
$ws1.Cells.Item(496,1).Value = "20:12:40"
$ws1.Cells.Item(496,2).Value = "22:07"
$ws1.Cells.Item(496,3).Value = "17_ROMERO"
$ws1.Cells.Item(496,4).Value = 115
$ws1.Cells.Item(497,2).Value = "22:07"
$ws1.Cells.Item(497,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(497,4).Value = 94
$ws1.Cells.Item(498,1).Value = "20:47:47"
$ws1.Cells.Item(498,2).Value = "22:08"
$ws1.Cells.Item(498,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(498,4).Value = 81
$ws1.Cells.Item(499,1).Value = "20:47:47"
$ws1.Cells.Item(499,2).Value = "22:13"
$ws1.Cells.Item(499,3).Value = "15_ABASTO"
$ws1.Cells.Item(499,4).Value = 86
$ws1.Cells.Item(499,5).Value = "LP1912"
$ws1.Cells.Item(500,1).Value = "20:33:52"
$ws1.Cells.Item(500,2).Value = "22:23"
$ws1.Cells.Item(500,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(500,4).Value = 110
$ws1.Cells.Item(500,5).Value = "LP1912"
$ws1.Cells.Item(501,1).Value = "20:33:52"
$ws1.Cells.Item(501,2).Value = "22:27"
$ws1.Cells.Item(501,3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(501,4).Value = 114
$ws1.Cells.Item(501,5).Value = "LP1912"
$ws1.Cells.Item(502,1).Value = "20:47:47"
$ws1.Cells.Item(502,2).Value = "22:27"
$ws1.Cells.Item(502,3).Value = "10_OLMOS"
$ws1.Cells.Item(502,4).Value = 100
$ws1.Cells.Item(502,5).Value = "LP1912"
$ws1.Cells.Item(503,1).Value = "20:33:52"
$ws1.Cells.Item(503,2).Value = "22:28"
$ws1.Cells.Item(503,3).Value = "10_OLMOS"
$ws1.Cells.Item(503,4).Value = 115
$ws1.Cells.Item(503,5).Value = "LP1912"
$ws1.Cells.Item(504,1).Value = "20:47:47"
$ws1.Cells.Item(504,2).Value = "22:39"
$ws1.Cells.Item(504,3).Value = "215A_EL PATO"
$ws1.Cells.Item(504,4).Value = 112
$ws1.Cells.Item(504,5).Value = "LP1912"

# ---- Sheet: LP1912-215 ----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2,1).Value = "Última actualización: 20:47:47"
$ws2.Cells.Item(3,1).Value = "Total filas: 51"
$ws2.Cells.Item(56,1).Value = "20:47:47"
$ws2.Cells.Item(56,2).Value = "22:39"
$ws2.Cells.Item(56,3).Value = "215A_EL PATO"
$ws2.Cells.Item(56,4).Value = 112
$ws2.Cells.Item(56,5).Value = "LP1912"

# ---- Sheet: 6203-6173 ----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2,1).Value = "Última actualización: 20:47:47"

